$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.033.28"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.416.35"
$ws.Range("D3").ClearFormats()
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.16"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.60"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.483"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.97"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.24%  "
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("E11").Value = "  +2.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.001.03"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.37"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.409.34"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000171"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.019.57"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.54%  "
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.52"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.99"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.65%  "
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.569"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.92"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.564.44"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("E26").Value = "  -2.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.179"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.62"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -2.79%  "
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  -1.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.05"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.73%  "
$ws.Range("E35").Value = "  +3.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.62"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.45%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "169.79"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.91"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "30.50"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.450.31"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0781"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.47"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("E44").Value = "  -1.86%  "
$ws.Range("E45").Value = "  -2.30%  "
$ws.Range("E46").Value = "  -3.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.543.19"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.91"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.82"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.20"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.81%  "
$ws.Range("E51").Value = "  +0.01%  "
